$d = $word.ActiveDocument

# The student's group label "3M" (stored as two adjacent runs "3" and "M")
# becomes "7N" (a single run).
$full = $d.Content.Text
$idx = $full.IndexOf("3M")

if ($idx -ge 0) {
    $r = $d.Range($idx, $idx + 2)
    $r.Text = "7N"
}

$d.Save()
